$wb = $excel.ActiveWorkbook

# --- Sheet 1: "All Published Values" -------------------------------------
$ws = $wb.Worksheets.Item("All Published Values")

# New row of published-rate data (row 11), mirroring the existing rows.
# Every column in this table is stored as text (even the numeric-looking
# ones), so force Text formatting on the destination range before typing
# the values in order to stop Excel's automatic number/date inference,
# then clear the formatting again so the new cells don't keep a style
# index that the sibling data rows don't have.
$newRow = $ws.Range("A11:J11")
$newRow.NumberFormat = "@"

$values = @(
    "2026-01-02",
    "2026-01-02 19:53:09",
    "697.85",
    "697.85",
    "700.79",
    "700.79",
    "702.88",
    "2026/01/02 19:53:09",
    "2026-01-02 11:57:40",
    "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(11, $i + 1).Value = $values[$i]
}

$newRow.ClearFormats()

# Re-establish the AutoFilter over the now-larger range (A1:J11).
$ws.AutoFilterMode = $false
$ws.Range("A1:J11").AutoFilter() | Out-Null

# Keep the sheet's hidden _FilterDatabase defined name in sync with the
# new AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$11"
    }
}

# --- Sheet 2: "Daily Summary" ---------------------------------------------
# The day's publish count increments from 9 to 10 with the new reading.
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Cells.Item(4, 2).Value = 10
